$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data as scraped on Mon Apr 29 08:37:58 UTC 2024
$ws.Range('D2').Value = '62.426.21'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '3.171.33'
$ws.Range('E3').Value = '  -4.52%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.44'
$ws.Range('E5').Value = '  -2.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.11'
$ws.Range('E6').Value = '  -6.43%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.169.13'
$ws.Range('E8').Value = '  -4.51%  '
$ws.Range('E9').Value = '  -4.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('E10').Value = '  -6.39%  '
$ws.Range('E11').Value = '  -6.64%  '
$ws.Range('E12').Value = '  -5.44%  '
$ws.Range('E13').Value = '  -6.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.90'
$ws.Range('E14').Value = '  -5.73%  '
$ws.Range('D15').Value = '3.704.90'
$ws.Range('E15').Value = '  -4.17%  '
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Value = '3.181.17'
$ws.Range('E17').Value = '  -4.25%  '
$ws.Range('D18').Value = '62.420.86'
$ws.Range('E18').Value = '  -2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.52'
$ws.Range('E19').Value = '  -5.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '454.82'
$ws.Range('E20').Value = '  -5.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.76'
$ws.Range('E21').Value = '  -3.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.699'
$ws.Range('E22').Value = '  -5.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.59'
$ws.Range('E23').Value = '  -5.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.30'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.98'
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.66'
$ws.Range('E28').Value = '  -4.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.89'
$ws.Range('E29').Value = '  -6.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.76'
$ws.Range('E30').Value = '  -5.34%  '
$ws.Range('E31').Value = '  -8.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.13'
$ws.Range('E32').Value = '  -7.97%  '
$ws.Range('E33').Value = '  -4.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.37'
$ws.Range('E34').Value = '  -7.43%  '
$ws.Range('E35').Value = '  -7.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.06'
$ws.Range('E37').Value = '  -3.44%  '
$ws.Range('D38').Value = '0.0₃0687'
$ws.Range('E38').Value = '  -9.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0381'
$ws.Range('E39').Value = '  -5.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '409.08'
$ws.Range('E40').Value = '  -5.47%  '
$ws.Range('D41').Value = '2.926.15'
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.112'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.94'
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.61'
$ws.Range('E44').Value = '  -5.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.248'
$ws.Range('E45').Value = '  -7.46%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.12'
$ws.Range('E47').Value = '  -4.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.42'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.72'
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.16'
$ws.Range('E50').Value = '  -5.28%  '
$ws.Range('E51').Value = '  -4.67%  '
